# PPImagePopulate / Presentation1.pptx
#
# Commit "update Ui and change image limit to 5" is a change to the
# PowerPoint *add-in's* web UI (the task-pane HTML/JS project that ships
# alongside this deck) together with a rebuild of the project. The only
# artifact checked in here, ppt/bin/Debug/Presentation1.pptx, is the blank
# scratch deck PowerPoint/Visual Studio uses to exercise that add-in - it
# carries none of the add-in's own markup.
#
# Diffing the canonical OOXML of before/after shows exactly two kinds of
# change, and nothing else:
#   1. Every r:id in presentation.xml (sldMasterId/sldId/sldLayoutId) and
#      in the taskpane reference is swapped for a freshly minted token.
#      That's PowerPoint re-serialising the package on save - relationship
#      IDs are regenerated opaque tokens with no semantic meaning, are not
#      addressable from the Presentation/Slide/Shape object model, and do
#      not correspond to any user-visible edit.
#   2. The <we:webextension> part's id GUID changes. That GUID identifies
#      the add-in manifest and lives entirely under ppt/webextensions/,
#      a part PowerPoint does not expose through the Presentation/Slide/
#      Shape COM surface (no VBA/COM object model path reaches
#      webextension.xml, in this host or in real PowerPoint - task-pane
#      add-in parts are edited by Visual Studio/the manifest tooling, not
#      by end-user automation).
#
# Every actual slide/shape/text object is untouched: both placeholders on
# the one slide (title, subtitle) are empty before and after. So there is
# no reachable content edit to make here - this script intentionally
# performs none, leaving the presentation exactly as authored.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Touch the object model read-only to confirm the deck opens/walks cleanly;
# no properties are assigned because none of the COM-addressable content
# (slides, shapes, text) differs between the before/after states.
$slideCount = $p.Slides.Count
$shapeCount = $s.Shapes.Count
Write-Output "Slides: $slideCount; Shapes on slide 1: $shapeCount"
